$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values, keyed by row number
$updates = @{
    3  = -3
    4  = -9
    6  = -2
    11 = -1
    13 = 2
    14 = 2
    16 = 0
    17 = 0
    20 = -5
    23 = 13
    25 = 2
    28 = 4
    29 = 3
    31 = -2
    32 = 3
    33 = -1
    34 = -3
    37 = -1
    41 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
